$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericText($ws, $addr, $val) {
    # Force a numeric-looking string to be stored as TEXT (not auto-converted
    # to a number) while keeping the cell's existing (default) style.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# --- TRAIN column (B) updates, written first so new shared strings are appended
#     in the same order the original authoring tool used ---
Set-NumericText $ws "B6" "0.7811080835603996"
Set-NumericText $ws "B7" "0.8138056312443234"
Set-NumericText $ws "B8" "0.046511627906976744"
Set-NumericText $ws "B9" "1.0"
$ws.Range("B10").Value = "[547, 276, 1097, 1113, 1185, 500, 198, 560, 1150, 1078, 564, 203, 1146, 572, 549, 283, 1220, 281, 1110, 1109, 1149, 544, 1288, 264, 239, 261, 257, 1080, 179, 1099, 1232, 1089, 217, 527, 293, 609, 260, 573, 1214, 1122, 1112, 1094, 534, 246, 1213, 237, 1170, 266, 577, 1176, 535, 1210, 505, 1086, 1123, 1138, 235, 324, 1286, 1305, 214, 556, 1115, 1313, 565, 1203, 271, 184, 279, 1100, 310, 312, 1230, 1178, 175, 294, 1111, 234, 196, 211, 562, 248, 301, 1129, 190, 288, 317, 531, 1200, 542, 302, 263, 1175, 1160, 499, 218, 1072, 1198, 1107, 1306, 222, 1209, 538, 290, 209, 200, 241, 579, 1294, 1307, 513, 1124, 278, 1312, 1166, 205, 602, 193, 1093, 553, 533, 526, 1295, 219, 578, 1190, 255, 1201, 1195, 557, 202, 1151, 1229, 1074, 1165, 1299, 585, 1130, 525, 1310, 307, 1205, 1298, 208, 215, 537, 524, 178, 187, 536, 1311, 1225, 292, 243, 199, 206, 603, 584, 1187, 509, 299, 501, 507, 575, 1087, 1134, 558, 1104, 1236, 240, 1145, 1147, 195, 519, 1131, 583, 305, 287, 571, 212, 566, 1105, 213, 256, 177, 599, 306, 582, 197, 1090, 245, 265, 1301, 1208, 554, 1168, 561, 181, 242, 309, 289, 270, 230, 286, 285]"

# --- TEST column (C) updates ---
$ws.Range("C5").Value = "TEST"
Set-NumericText $ws "C6" "0.7709090909090909"
Set-NumericText $ws "C7" "0.78"
Set-NumericText $ws "C8" "0.050980392156862744"
Set-NumericText $ws "C9" "1.0"
$ws.Range("C10").Value = "[517, 1101, 504, 1083, 315, 232, 1084, 182, 1308, 601, 1081, 1119, 1182, 1096, 227, 233, 1157, 1291, 1211, 229, 1140, 1116, 1224, 1290, 1079, 204, 250, 503, 520, 495, 580, 1289, 1073, 1164, 521, 221, 1117, 522, 1127, 1143, 1207, 224, 1219, 551, 295, 1300, 1095, 236, 1216, 606, 529, 252, 545, 1085, 555, 1231, 1139, 201, 1125, 552, 1303, 1133, 313, 1173, 1171, 543, 550, 1215, 559, 180, 1152, 1128, 1155, 297, 1304, 514, 1142, 512, 207, 210, 581, 1194, 532, 1227, 249, 192, 511, 280, 1235, 1234, 311, 1106, 1144, 1154, 296, 282, 300, 1218, 1075, 268, 1186, 1126, 1181, 259, 1223, 1191, 191, 1141, 1172, 1233, 1189, 238, 1102, 277, 1082, 1206, 1297, 548, 516, 605, 183, 273, 518, 598, 1197, 1156, 506, 604, 244, 189, 567, 496, 1184, 1221, 228, 1135, 1302, 574, 502, 274, 563, 1180, 1314, 272, 194, 493, 1179, 225, 1161, 1204, 262, 508, 1177, 1108, 494, 608, 254, 1315, 318, 1136, 530, 1193, 1292, 291, 220, 1092, 607, 267, 1174, 1158, 1226, 298, 304, 1285, 314, 1169, 597, 1121, 223, 1163, 1088, 1202, 528, 253, 570, 185, 541, 568, 1077, 569, 1114, 1103, 1098, 576, 188, 1296, 1132, 258, 1120, 1287, 176, 186, 251, 1148, 1188, 1309, 1199, 1137, 275, 269, 1293, 1167, 308, 498, 1162, 284, 1159, 540, 231, 497, 1183, 1091, 1118, 316, 523, 1222, 515, 510, 226, 303, 600, 216, 1228, 539, 546, 1076, 247, 1217, 1153, 1196, 1212, 1192]"

# --- BNN block relabelling (rows shift up after the layer-fidelity rows are removed) ---
$ws.Range("A13").Value = "BNN"
$ws.Range("A14").Value = "# ENTRIES:"
$ws.Range("A18").Value = "# AVERAGE THRESHOLDS PER USED NEURON:"

# --- Remove the trailing "CONDITIONS/FIDELITY PER LAYER" rows (19-21) entirely ---
$ws.Range("A19:B21").ClearContents()
